$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: correct effort split between "Effort" and "Additional Effort"
$ws.Range("B8").Value = 2
$ws.Range("C8").Value = 2

# Row 9: new entry for mutex implementation work
$ws.Range("A9").Value = $ws.Range("A8").Value2
$ws.Range("B9").Value = 2.25
$ws.Range("D9").Value = "Implementation of mutexes"

# Row 10: new entry with detailed status note
$ws.Range("A10").Value = 41435
$ws.Range("B10").Value = 2
$ws.Range("D10").Value = "Implementation of mutexes: Basically done. No test case implemented yet, no testing done yet"

# Move the active selection to reflect where the user ended up editing
$ws.Range("E10").Select()
